# TC_63769_Verify_Current(DC_Units)_Calculation_FC.xlsx
# Adds a "Loop"/"Column" lookup block (J1:K5) to both worksheets, updates
# the "Current (DC Units)" values (column G) on each sheet, and leaves
# "Update Devices" as the active/selected sheet (matching the saved view
# state captured in the target workbook).

$wb = $excel.ActiveWorkbook

$wsLoopA = $wb.Worksheets.Item("Add Devices Loop A")
$wsUpdate = $wb.Worksheets.Item("Update Devices")

# ---------------------------------------------------------------------
# Sheet "Add Devices Loop A"
# ---------------------------------------------------------------------

# Updated "Current (DC Units)" reading.
$wsLoopA.Range("G1").Value = 313

# New "Loop" / "Column" header cells - copy the format already used by the
# neighbouring "DC Unit Loading Details Name" header (E1) then set text.
$wsLoopA.Range("E1").Copy()
$wsLoopA.Range("J1:K1").PasteSpecial(-4122)
$wsLoopA.Range("J1").Value = "Loop"
$wsLoopA.Range("K1").Value = "Column"

# New "Built-in Loop-A".."Built-in Loop-D" rows - copy the plain bordered
# body format already used lower in the same table (B4) then set text.
$wsLoopA.Range("B4").Copy()
$wsLoopA.Range("J2:K5").PasteSpecial(-4122)
$wsLoopA.Range("J2").Value = "Built-in Loop-A"
$wsLoopA.Range("J3").Value = "Built-in Loop-B"
$wsLoopA.Range("J4").Value = "Built-in Loop-C"
$wsLoopA.Range("J5").Value = "Built-in Loop-D"
$wsLoopA.Range("K2").Value = 2
$wsLoopA.Range("K3").Style = "Normal"
$wsLoopA.Range("K4").Style = "Normal"
$wsLoopA.Range("K5").Style = "Normal"

$wsLoopA.Activate()
$excel.ActiveWindow.ScrollColumn = 6
$wsLoopA.Range("J1:K5").Select()

# ---------------------------------------------------------------------
# Sheet "Update Devices"
# ---------------------------------------------------------------------

# Updated "Current (DC Units)" readings.
$wsUpdate.Range("G2").Value = 339
$wsUpdate.Range("G3").Value = 336
$wsUpdate.Range("G4").Value = 311

# New "Loop" / "Column" header cells.
$wsUpdate.Range("E1").Copy()
$wsUpdate.Range("J1:K1").PasteSpecial(-4122)
$wsUpdate.Range("J1").Value = "Loop"
$wsUpdate.Range("K1").Value = "Column"

# New "Built-in Loop-A".."Built-in Loop-D" rows (column K on this sheet
# only carries the numeric index on row 2, rows 3-5 stay blank/removed).
$wsUpdate.Range("B4").Copy()
$wsUpdate.Range("J2:J5").PasteSpecial(-4122)
$wsUpdate.Range("J2").Value = "Built-in Loop-A"
$wsUpdate.Range("J3").Value = "Built-in Loop-B"
$wsUpdate.Range("J4").Value = "Built-in Loop-C"
$wsUpdate.Range("J5").Value = "Built-in Loop-D"

$wsUpdate.Range("B4").Copy()
$wsUpdate.Range("K2").PasteSpecial(-4122)
$wsUpdate.Range("K2").Value = 2

# The old placeholder cells K3:K5 are dropped entirely on this sheet.
$wsUpdate.Range("K3").Clear()
$wsUpdate.Range("K4").Clear()
$wsUpdate.Range("K5").Clear()

# "Update Devices" ends up the active sheet/tab with the new block selected.
$wsUpdate.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$wsUpdate.Range("J1:K5").Select()
